# Fill in the previously-empty execution-time table cells (and fix up one
# existing value) in the "Execution Time for Different Clone Detection
# Tools" table. Cell(row, col) is 1-based: row 1 is the header row, and
# columns are 1=tool, 2=RQ, 3=Ctags, 4=BrlCad, 5=Freecol, 6=Carol, 7=Jabref.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-CellSimpleText($row, $col, $text) {
    $table = $d.Tables(1)
    $cell = $table.Cell($row, $col)
    $xml = "<w:p $wNs><w:r><w:t>$text</w:t></w:r></w:p>"
    $cell.Range.InsertXML($xml)
}

function Set-CellRunsXml($row, $col, $innerRunsXml) {
    $table = $d.Tables(1)
    $cell = $table.Cell($row, $col)
    $xml = "<w:p $wNs>$innerRunsXml</w:p>"
    $cell.Range.InsertXML($xml)
}

# --- NiCad / RQ3 ---
Set-CellSimpleText 4 4 "2 m 10 s"

# --- Deckard (2nd) / RQ1 ---
$runs = "<w:r><w:t>2</w:t></w:r>" +
        "<w:r><w:t>6</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> m 4</w:t></w:r>" +
        "<w:r><w:t>5</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> s</w:t></w:r>"
Set-CellRunsXml 10 3 $runs
Set-CellSimpleText 10 4 "38 m 5 s"
Set-CellSimpleText 10 5 "82 h 37 m"
Set-CellSimpleText 10 6 "17 h 55 m"

# --- Deckard (2nd) / RQ2 ---
$runs = "<w:r><w:t>2</w:t></w:r>" +
        "<w:r><w:t>9</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> m </w:t></w:r>" +
        "<w:r><w:t>29</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> s</w:t></w:r>"
Set-CellRunsXml 11 3 $runs
Set-CellSimpleText 11 4 "47 m 30 s"
Set-CellSimpleText 11 7 "19 h 21 m"

# --- Deckard (2nd) / RQ3 ---
Set-CellSimpleText 12 3 "25 m 59 s"
Set-CellSimpleText 12 4 "41 m 29 s"
Set-CellSimpleText 12 7 "17 h 59 m"

# --- Deckard (2nd) / RQ4 ---
Set-CellSimpleText 13 3 "12 m 37 s"
Set-CellSimpleText 13 4 "5 m 15 s"
Set-CellSimpleText 13 7 "3 h 31 m"

# --- ConQat / RQ1 ---
Set-CellSimpleText 14 3 "4 h 19 m"
Set-CellSimpleText 14 4 "2 h 33 m"

# --- ConQat / RQ2 ---
Set-CellSimpleText 15 3 "8 h 17 m"
Set-CellSimpleText 15 4 "4 h 2 m"

# --- ConQat / RQ3 ---
Set-CellSimpleText 16 3 "8 h 26 m"
Set-CellSimpleText 16 4 "3 h 57 m"

Write-Output "Table cells updated."
